$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $maxRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
